$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comp controls")

$ws.Range("B4").Value  = "LIVE GREEN:B515"
$ws.Range("B5").Value  = "CD197:G560"
$ws.Range("B6").Value  = "CD4:B710"
$ws.Range("B7").Value  = "CD45RA:G780"
$ws.Range("B8").Value  = "CD194:G780"
$ws.Range("B9").Value  = "CD27:G780"
$ws.Range("B10").Value = "CD11c:G780"
$ws.Range("B11").Value = "CD196:G780"
$ws.Range("B12").Value = "CD38:R660"
$ws.Range("B13").Value = "CD127:R660"
$ws.Range("B14").Value = "CD8:R780"
$ws.Range("B15").Value = "CD45RO:R780"
$ws.Range("B16").Value = "CD20:R780"
$ws.Range("B17").Value = "CD3+19+20:R780"
$ws.Range("B18").Value = "CD3:V450"
$ws.Range("B19").Value = "HLA-DR:V545"
